$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.722.99'
$ws.Range('E2').Value = '  +0.42%  '
$ws.Range('D3').Value = '3.099.13'
$ws.Range('E3').Value = '  +3.86%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '''388.62'
$ws.Range('E5').Value = '  +1.90%  '
$ws.Range('D6').Value = '''103.53'
$ws.Range('E6').Value = '  -0.45%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  +0.11%  '
$ws.Range('D9').Value = '''0.588'
$ws.Range('E9').Value = '  -1.05%  '
$ws.Range('E10').Value = '  +0.92%  '
$ws.Range('E11').Value = '  -0.07%  '
$ws.Range('D12').Value = '''0.0863'
$ws.Range('E12').Value = '  +0.55%  '
$ws.Range('D13').Value = '3.584.83'
$ws.Range('E13').Value = '  +3.77%  '
$ws.Range('D14').Value = '''18.66'
$ws.Range('E14').Value = '  +0.92%  '
$ws.Range('E15').Value = '  -0.84%  '
$ws.Range('D16').Value = '3.098.89'
$ws.Range('E16').Value = '  +3.85%  '
$ws.Range('D17').Value = '''0.982'
$ws.Range('E17').Value = '  -1.58%  '
$ws.Range('E18').Value = '  -5.15%  '
$ws.Range('D19').Value = '51.880.32'
$ws.Range('E19').Value = '  +0.71%  '
$ws.Range('D20').Value = '''3.19'
$ws.Range('E20').Value = '  +3.04%  '
$ws.Range('D21').Value = '''12.47'
$ws.Range('E21').Value = '  -0.89%  '
$ws.Range('D22').Value = '0.0₃0970'
$ws.Range('E22').Value = '  +0.76%  '
$ws.Range('D23').Value = '''70.04'
$ws.Range('E23').Value = '  -0.39%  '
$ws.Range('D24').Value = '''268.75'
$ws.Range('E24').Value = '  +0.49%  '
$ws.Range('D25').Value = '''3.12'
$ws.Range('E25').Value = '  -3.21%  '
$ws.Range('E26').Value = '  +3.30%  '
$ws.Range('D27').Value = '''27.06'
$ws.Range('E27').Value = '  +3.62%  '
$ws.Range('E28').Value = '  +0.74%  '
$ws.Range('D29').Value = '''7.21'
$ws.Range('E29').Value = '  -0.95%  '
$ws.Range('E30').Value = '  +0.07%  '
$ws.Range('E31').Value = '  -0.13%  '
$ws.Range('D32').Value = '''10.35'
$ws.Range('E32').Value = '  -0.76%  '
$ws.Range('D33').Value = '''35.64'
$ws.Range('E33').Value = '  +2.93%  '
$ws.Range('E34').Value = '  +0.68%  '
$ws.Range('D35').Value = '''50.44'
$ws.Range('E35').Value = '  -1.81%  '
$ws.Range('E36').Value = '  +0.32%  '
$ws.Range('E37').Value = '  -0.14%  '
$ws.Range('D38').Value = '''3.41'
$ws.Range('E38').Value = '  +3.84%  '
$ws.Range('D39').Value = '''0.291'
$ws.Range('E39').Value = '  +7.30%  '
$ws.Range('E40').Value = '  +1.84%  '
$ws.Range('E41').Value = '  +0.05%  '
$ws.Range('E42').Value = '  +0.70%  '
$ws.Range('B43').Value = 'Monero'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D43').Value = '''127.52'
$ws.Range('E43').Value = '  +0.81%  '
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D44').Value = '''0.116'
$ws.Range('E44').Value = '  -0.60%  '
$ws.Range('E45').Value = '  -3.85%  '
$ws.Range('D46').Value = '''22.20'
$ws.Range('E46').Value = '  +3.49%  '
$ws.Range('D47').Value = '''2.48'
$ws.Range('E47').Value = '  +5.15%  '
$ws.Range('E48').Value = '  +2.09%  '
$ws.Range('D49').Value = '2.048.01'
$ws.Range('E49').Value = '  +1.14%  '
$ws.Range('D50').Value = '3.405.13'
$ws.Range('E50').Value = '  +3.89%  '
$ws.Range('E51').Value = '  +6.10%  '
